$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.476561903953552
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 2.410345315933228
$ws.Range("D1").Value = 1.562836408615112
$ws.Range("E1").Value = 1.275827646255493
